$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111790785
$ws.Range("B2").Value = 77515
$ws.Range("C2").Value = "Ovaliderad"
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = "Garnlav"
$ws.Range("G2").Value = "Alectoria sarmentosa"
$ws.Range("H2").Value = "(Ach.) Ach."
$ws.Range("P2").Value = "Åsvarptjärnen (Åsvarptjärnen), Jmt"
$ws.Range("Q2").Value = 489818.2822038208
$ws.Range("R2").Value = 6949032.207674611
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = "Jämtland"
$ws.Range("U2").Value = "Berg"
$ws.Range("V2").Value = "Jämtland"
$ws.Range("W2").Value = "Hackås"
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2023-08-30"
$ws.Range("Z2").Value = "18:34"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2023-08-30"
$ws.Range("AB2").Value = "18:34"
$ws.Range("AC2").Value = "Många träd med mycket lav i området"
$ws.Range("AD2").Value = $false
$ws.Range("AE2").Value = $false
$ws.Range("AG2").Value = $false
$ws.Range("AW2").Value = "Erik Wilhelmsson"
$ws.Range("AX2").Value = "Erik Wilhelmsson"

# Row 3
$ws.Range("A3").Value = 111790625
$ws.Range("B3").Value = 96348
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = "Knärot"
$ws.Range("G3").Value = "Goodyera repens"
$ws.Range("H3").Value = "(L.) R. Br."
$ws.Range("P3").Value = "Åsvarptjärnen (Åsvarptjärnen), Jmt"
$ws.Range("Q3").Value = 489824.6884970492
$ws.Range("R3").Value = 6949020.70113107
$ws.Range("S3").Value = 1
$ws.Range("T3").Value = "Jämtland"
$ws.Range("U3").Value = "Berg"
$ws.Range("V3").Value = "Jämtland"
$ws.Range("W3").Value = "Hackås"
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2023-08-30"
$ws.Range("Z3").Value = "18:29"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2023-08-30"
$ws.Range("AB3").Value = "18:29"
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AW3").Value = "Erik Wilhelmsson"
$ws.Range("AX3").Value = "Erik Wilhelmsson"

# Row 4
$ws.Range("A4").Value = 111792337
$ws.Range("B4").Value = 96348
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "VU"
$ws.Range("E4").Value = 220787
$ws.Range("F4").Value = "Knärot"
$ws.Range("G4").Value = "Goodyera repens"
$ws.Range("H4").Value = "(L.) R. Br."
$ws.Range("P4").Value = "Åsvarptjärnen (Åsvarptjärnen), Jmt"
$ws.Range("Q4").Value = 489763.7116335144
$ws.Range("R4").Value = 6949091.647604217
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = "Jämtland"
$ws.Range("U4").Value = "Berg"
$ws.Range("V4").Value = "Jämtland"
$ws.Range("W4").Value = "Hackås"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2023-08-30"
$ws.Range("Z4").Value = "19:22"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2023-08-30"
$ws.Range("AB4").Value = "19:22"
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AW4").Value = "Erik Wilhelmsson"
$ws.Range("AX4").Value = "Erik Wilhelmsson"

# Row 5
$ws.Range("A5").Value = 111919588
$ws.Range("B5").Value = 78578
$ws.Range("C5").Value = "Ovaliderad"
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 6458
$ws.Range("F5").Value = "Lunglav"
$ws.Range("G5").Value = "Lobaria pulmonaria"
$ws.Range("H5").Value = "(L.) Hoffm."
$ws.Range("P5").Value = "Gillhov, Jmt"
$ws.Range("Q5").Value = 490132.6970436619
$ws.Range("R5").Value = 6948774.399184751
$ws.Range("S5").Value = 1
$ws.Range("T5").Value = "Jämtland"
$ws.Range("U5").Value = "Berg"
$ws.Range("V5").Value = "Jämtland"
$ws.Range("W5").Value = "Hackås"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2023-09-05"
$ws.Range("Z5").Value = "00:00"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2023-09-05"
$ws.Range("AB5").Value = "00:00"
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AW5").Value = "Erik Wilhelmsson"
$ws.Range("AX5").Value = "Erik Wilhelmsson"

# Row 6
$ws.Range("A6").Value = 111915405
$ws.Range("B6").Value = 88946
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "VU"
$ws.Range("E6").Value = 256335
$ws.Range("F6").Value = "Taggfingersvamp"
$ws.Range("G6").Value = "Ramaria karstenii"
$ws.Range("H6").Value = "(Sacc. & P.Syd.) Corner"
$ws.Range("P6").Value = "Nordvallen (Nordvallen), Jmt"
$ws.Range("Q6").Value = 490124.2549094426
$ws.Range("R6").Value = 6948875.054878937
$ws.Range("S6").Value = 1
$ws.Range("T6").Value = "Jämtland"
$ws.Range("U6").Value = "Berg"
$ws.Range("V6").Value = "Jämtland"
$ws.Range("W6").Value = "Hackås"
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2023-09-05"
$ws.Range("Z6").Value = "00:00"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2023-09-05"
$ws.Range("AB6").Value = "00:00"
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AW6").Value = "Erik Wilhelmsson"
$ws.Range("AX6").Value = "Erik Wilhelmsson"

# Clear AC3/AC4 since the public-comment text moved to AC2 in the reordered rows
$ws.Range("AC3").ClearContents()
$ws.Range("AC4").ClearContents()